# Write a df into a specific position of the excel file.
#
# The original sheet holds a small "dataframe-like" block at A1:D3
# (header row "x,1,2,3", a "y" row of same-column references, and a "z"
# row that multiplies the header by the "y" row). This change writes the
# same kind of block again, offset by one row/column (anchored at B2),
# overlapping the bottom-right of the original block and extending the
# sheet to E4.
#
# Concretely:
#   - E2 gains the literal (text) value "3"
#   - B3/C3/D3's multiplication formulas now multiply the header by the
#     column to their *left* (row 2) instead of their own column
#   - E3 becomes a simple same-row reference to E2 (like the "y" row
#     pattern)
#   - a new row 4 is added: B4 is the text label "z", and C4/D4/E4 repeat
#     the "multiply by column to the left" formula pattern one row down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 = "3" as TEXT (not the number 3) -- force text via NumberFormat,
# then restore the default "Normal" style so no stray number-format
# sticks around on the cell.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3"
$ws.Range("E2").Style = "Normal"

# Row 3: rewrite the multiplication formulas to reference the previous
# column's row-2 cell instead of their own column's row-2 cell. The
# left-most one (B3) loses its second operand entirely, matching the
# source data exactly.
$ws.Range("B3").Formula = "=(B1 * )"
$ws.Range("C3").Formula = "=(C1 * B2)"
$ws.Range("D3").Formula = "=(D1 * C2)"
$ws.Range("E3").Formula = "=E2"

# New row 4.
$ws.Range("B4").Value = "z"
$ws.Range("C4").Formula = "=(C2 * )"
$ws.Range("D4").Formula = "=(D2 * C3)"
$ws.Range("E4").Formula = "=(E2 * D3)"
